# Auto-applied update from scheduled market-data runner.
# Updates crafting-profit value cells (currentAveragePrice / LevePrice / LeveProfit columns)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match refreshed market data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 482.9
$ws.Range("I15").Value = 482.9
$ws.Range("K15").Value = 1448.7
$ws.Range("M15").Value = -1279.7
$ws.Range("H33").Value = 259.72726
$ws.Range("I33").Value = 266.5862
$ws.Range("J33").Value = 210
$ws.Range("K33").Value = 266.5862
$ws.Range("L33").Value = 210
$ws.Range("M33").Value = -37.58620000000002
$ws.Range("N33").Value = -668
$ws.Range("H113").Value = 2801.6
$ws.Range("I113").Value = 1863
$ws.Range("J113").Value = 3427.3333
$ws.Range("K113").Value = 1863
$ws.Range("L113").Value = 3427.3333
$ws.Range("M113").Value = 1391
$ws.Range("N113").Value = -9935.3333

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7229
$ws.Range("I61").Value = 5464.25
$ws.Range("J61").Value = 11262.714
$ws.Range("K61").Value = 5464.25
$ws.Range("L61").Value = 11262.714
$ws.Range("M61").Value = -5252.25
$ws.Range("N61").Value = -11686.714
$ws.Range("H63").Value = 1540.6666
$ws.Range("I63").Value = 1394.2858
$ws.Range("J63").Value = 2053
$ws.Range("K63").Value = 1394.2858
$ws.Range("L63").Value = 2053
$ws.Range("M63").Value = -708.2858000000001
$ws.Range("N63").Value = -3425
$ws.Range("H66").Value = 1540.6666
$ws.Range("I66").Value = 1394.2858
$ws.Range("J66").Value = 2053
$ws.Range("K66").Value = 6971.429
$ws.Range("L66").Value = 10265
$ws.Range("M66").Value = -3539.429
$ws.Range("N66").Value = -17129
$ws.Range("H136").Value = 7229
$ws.Range("I136").Value = 5464.25
$ws.Range("J136").Value = 11262.714
$ws.Range("K136").Value = 16392.75
$ws.Range("L136").Value = 33788.142
$ws.Range("M136").Value = -13842.75
$ws.Range("N136").Value = -38888.142

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 9800
$ws.Range("J6").Value = 9800
$ws.Range("L6").Value = 9800
$ws.Range("N6").Value = -10026
$ws.Range("H51").Value = 49961
$ws.Range("J51").Value = 49961
$ws.Range("L51").Value = 49961
$ws.Range("N51").Value = -50943
$ws.Range("H107").Value = 2462.6843
$ws.Range("J107").Value = 2639.125
$ws.Range("L107").Value = 2639.125
$ws.Range("N107").Value = -6479.125
$ws.Range("H134").Value = 1900.3112
$ws.Range("I134").Value = 1845.5714
$ws.Range("K134").Value = 5536.7142
$ws.Range("M134").Value = -3001.7142
$ws.Range("H140").Value = 38378.5
$ws.Range("J140").Value = 38378.5
$ws.Range("L140").Value = 38378.5
$ws.Range("N140").Value = -48738.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100.5
$ws.Range("I16").Value = 960.25
$ws.Range("J16").Value = 1240.75
$ws.Range("K16").Value = 960.25
$ws.Range("L16").Value = 1240.75
$ws.Range("M16").Value = -673.25
$ws.Range("N16").Value = -1814.75
$ws.Range("H44").Value = 10064
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("H76").Value = 9846.154
$ws.Range("I76").Value = 9846.154
$ws.Range("K76").Value = 9846.154
$ws.Range("M76").Value = -9531.154
$ws.Range("H79").Value = 9846.154
$ws.Range("I79").Value = 9846.154
$ws.Range("K79").Value = 9846.154
$ws.Range("M79").Value = -8754.154
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("H105").Value = 1100
$ws.Range("I105").Value = 798.1667
$ws.Range("J105").Value = 2005.5
$ws.Range("K105").Value = 798.1667
$ws.Range("L105").Value = 2005.5
$ws.Range("M105").Value = 948.8333
$ws.Range("N105").Value = -5499.5
$ws.Range("H107").Value = 785.7917
$ws.Range("I107").Value = 844.4286
$ws.Range("K107").Value = 844.4286
$ws.Range("M107").Value = 1075.5714
$ws.Range("H113").Value = 1100.5
$ws.Range("I113").Value = 960.25
$ws.Range("J113").Value = 1240.75
$ws.Range("K113").Value = 960.25
$ws.Range("L113").Value = 1240.75
$ws.Range("M113").Value = 1209.75
$ws.Range("N113").Value = -5580.75
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 2653.16
$ws.Range("I132").Value = 2465.7693
$ws.Range("J132").Value = 2856.1667
$ws.Range("K132").Value = 7397.3079
$ws.Range("L132").Value = 8568.500100000001
$ws.Range("M132").Value = -4867.3079
$ws.Range("N132").Value = -13628.5001
$ws.Range("N44").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("N125").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 23815250
$ws.Range("I5").Value = 736.63635
$ws.Range("J5").Value = 50011216
$ws.Range("K5").Value = 2209.90905
$ws.Range("L5").Value = 150033648
$ws.Range("M5").Value = -2097.90905
$ws.Range("N5").Value = -150033872
$ws.Range("H20").Value = 1999.1666
$ws.Range("I20").Value = 900
$ws.Range("J20").Value = 2784.2856
$ws.Range("K20").Value = 2700
$ws.Range("L20").Value = 8352.856800000001
$ws.Range("M20").Value = -2473
$ws.Range("N20").Value = -8806.856800000001
$ws.Range("H40").Value = 70.63636
$ws.Range("I40").Value = 60.875
$ws.Range("K40").Value = 243.5
$ws.Range("M40").Value = -174.5
$ws.Range("H113").Value = 779.98865
$ws.Range("I113").Value = 786.9231
$ws.Range("J113").Value = 725.9
$ws.Range("K113").Value = 2360.7693
$ws.Range("L113").Value = 2177.7
$ws.Range("M113").Value = -190.7692999999999
$ws.Range("N113").Value = -6517.7
$ws.Range("H116").Value = 1504.8334
$ws.Range("I116").Value = 757.25
$ws.Range("K116").Value = 2271.75
$ws.Range("M116").Value = 1170.25
$ws.Range("H135").Value = 23815250
$ws.Range("I135").Value = 736.63635
$ws.Range("J135").Value = 50011216
$ws.Range("K135").Value = 6629.72715
$ws.Range("L135").Value = 450100944
$ws.Range("M135").Value = -4094.72715
$ws.Range("N135").Value = -450106014

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1997.9524
$ws.Range("I113").Value = 2003.8889
$ws.Range("J113").Value = 1993.5
$ws.Range("K113").Value = 2003.8889
$ws.Range("L113").Value = 1993.5
$ws.Range("M113").Value = 166.1111000000001
$ws.Range("N113").Value = -6333.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 685
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 556.25
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 556.25
$ws.Range("M22").Value = -905
$ws.Range("N22").Value = -1146.25
$ws.Range("H27").Value = 685
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 556.25
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 556.25
$ws.Range("M27").Value = -1093
$ws.Range("N27").Value = -770.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 39900
$ws.Range("J128").Value = 39900
$ws.Range("L128").Value = 39900
$ws.Range("N128").Value = -49860
$ws.Range("H132").Value = 1523.1482
$ws.Range("I132").Value = 1261.75
$ws.Range("J132").Value = 2045.9445
$ws.Range("K132").Value = 3785.25
$ws.Range("L132").Value = 6137.833500000001
$ws.Range("M132").Value = -1255.25
$ws.Range("N132").Value = -11197.8335
